$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the reference date (column G) from 2024-08-22 (45526) to 2024-08-23 (45527) for every data row
for ($r = 2; $r -le 274; $r++) {
    $ws.Cells.Item($r, 7).Value = 45527
}

# Refresh the "Saldo Previsto" (E) and "Vl. Total" (H) figures that came with the new daily extract
$ws.Cells.Item(5, 5).Value = 18922.68
$ws.Cells.Item(5, 8).Value = 18922.68
$ws.Cells.Item(8, 5).Value = 9914.2099999999991
$ws.Cells.Item(8, 8).Value = 9914.2099999999991
$ws.Cells.Item(15, 5).Value = 33915.89
$ws.Cells.Item(15, 8).Value = 33915.89
$ws.Cells.Item(17, 5).Value = 12982.64
$ws.Cells.Item(17, 8).Value = 12982.64
$ws.Cells.Item(43, 5).Value = 1047.21
$ws.Cells.Item(43, 8).Value = 1047.21
$ws.Cells.Item(49, 5).Value = 6937.61
$ws.Cells.Item(49, 8).Value = 6937.61
$ws.Cells.Item(52, 5).Value = 868.19
$ws.Cells.Item(52, 8).Value = 868.19
$ws.Cells.Item(60, 5).Value = 16991.8
$ws.Cells.Item(60, 8).Value = 16991.8
$ws.Cells.Item(99, 5).Value = 1054.3
$ws.Cells.Item(99, 8).Value = 1054.3
$ws.Cells.Item(104, 5).Value = 1186.8900000000001
$ws.Cells.Item(104, 8).Value = 1186.8900000000001
$ws.Cells.Item(108, 5).Value = 1189.01
$ws.Cells.Item(108, 8).Value = 1189.01
$ws.Cells.Item(112, 5).Value = 17000.38
$ws.Cells.Item(112, 8).Value = 17000.38
$ws.Cells.Item(120, 5).Value = 37409.360000000001
$ws.Cells.Item(120, 8).Value = 37409.360000000001
$ws.Cells.Item(132, 5).Value = 1024.1400000000001
$ws.Cells.Item(132, 8).Value = 1024.1400000000001
$ws.Cells.Item(143, 5).Value = 45126
$ws.Cells.Item(143, 8).Value = 45126
$ws.Cells.Item(158, 5).Value = 88.78
$ws.Cells.Item(158, 8).Value = 88.78
$ws.Cells.Item(173, 5).Value = 1112.3499999999999
$ws.Cells.Item(173, 8).Value = 1112.3499999999999
$ws.Cells.Item(235, 5).Value = 1041.49
$ws.Cells.Item(235, 8).Value = 1041.49
$ws.Cells.Item(249, 5).Value = 599.02
$ws.Cells.Item(249, 8).Value = 599.02
$ws.Cells.Item(264, 5).Value = 1138.8
$ws.Cells.Item(264, 8).Value = 1138.8
$ws.Cells.Item(265, 5).Value = 1074.42
$ws.Cells.Item(265, 8).Value = 1074.42
$ws.Cells.Item(270, 5).Value = 14871.79
$ws.Cells.Item(270, 8).Value = 14871.79
$ws.Cells.Item(271, 5).Value = 19253.349999999999
$ws.Cells.Item(271, 8).Value = 19253.349999999999
$ws.Cells.Item(273, 5).Value = 1046.6500000000001
$ws.Cells.Item(273, 8).Value = 1046.6500000000001

# Rename the sheet to match the new extraction timestamp
$ws.Name = "IClientBalance-20240823-115142-"

# Restore the active selection to C6
$ws.Range("C6").Select()
